# Weekly price-sheet update: a new week's record is inserted as row 22
# (Región de Arica y Parinacota, 2022-01-06, 50 x $/malla 20 kilos @ 18000),
# pushing the previously existing rows 22-37 down to rows 23-38.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 22; everything below (old rows 22-37)
# shifts down to rows 23-38, and the sheet dimension grows to A1:R38.
$ws.Rows("22:22").Insert()

# Populate the newly inserted row 22 with this week's data.
$ws.Range("A22").Value = 10
$ws.Range("B22").Value = "Vega Modelo de Temuco"
$ws.Range("C22").Value = "La Araucanía"
$ws.Range("D22").Value = 44567
$ws.Range("E22").Value = 9
$ws.Range("F22").Value = 100114002
$ws.Range("G22").Value = "Camote"
$ws.Range("H22").Value = "Sin especificar"
$ws.Range("I22").Value = "Primera"
$ws.Range("J22").Value = 50
$ws.Range("K22").Value = 18000
$ws.Range("L22").Value = 18000
$ws.Range("M22").Value = 18000
$ws.Range("N22").Value = "`$/malla 20 kilos"
$ws.Range("O22").Value = "Región de Arica y Parinacota"
$ws.Range("P22").Value = 900
$ws.Range("Q22").Value = 20
$ws.Range("R22").Value = "Hortaliza"
